# Fixed a bug in WinWeight
# Rewrites the data rows (2-23) of Sheet1 with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(201, 9, 30, 15, 45, 30)
    3  = @(1201, 2, 10, 10, 10, 10)
    4  = @(1202, 2, 10, 10, 10, 10)
    5  = @(1203, 3, 15, 15, 15, 15)
    6  = @(101, 9, 30, 15, 60, 15)
    7  = @(901, 16, 15, 45, 60, 60)
    8  = @(501, 9, 52, 30, 75, 45)
    9  = @(401, 9, 48, 67, 75, 45)
    10 = @(601, 9, 60, 67, 60, 42)
    11 = @(801, 3, 67, 65, 52, 45)
    12 = @(301, 6, 45, 30, 60, 45)
    14 = @(1001, 18, 30, 75, 60, 72)
    15 = @(701, 3, 90, 45, 97, 15)
    17 = @(1101, 0, 15, 30, 30, 0)
    18 = @(2, 0, 2, 2, 2, 2)
    19 = @(502, 0, 4, 0, 0, 0)
    20 = @(1, 0, 2, 2, 2, 2)
    22 = @(602, 0, 0, 4, 0, 9)
    23 = @(402, 0, 0, 4, 0, 0)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
